$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '35.522.11'
$r.Style = $s
$ws.Range('E2').Value = '  -2.54%  '
$r = $ws.Range('D3')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '1.974.37'
$r.Style = $s
$ws.Range('E3').Value = '  -3.88%  '
$ws.Range('E4').Value = '  +0.04%  '
$r = $ws.Range('D5')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '244.67'
$r.Style = $s
$ws.Range('E5').Value = '  +1.21%  '
$r = $ws.Range('D6')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.635'
$r.Style = $s
$ws.Range('E6').Value = '  -4.39%  '
$r = $ws.Range('D7')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '57.03'
$r.Style = $s
$ws.Range('E7').Value = '  +4.94%  '
$ws.Range('E8').Value = '  +0.05%  '
$r = $ws.Range('D9')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '58.62'
$r.Style = $s
$ws.Range('E9').Value = '  +0.52%  '
$r = $ws.Range('D10')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.358'
$r.Style = $s
$ws.Range('E10').Value = '  +0.67%  '
$r = $ws.Range('D11')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.0732'
$r.Style = $s
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('E12').Value = '  -3.02%  '
$r = $ws.Range('D13')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.949'
$r.Style = $s
$ws.Range('E13').Value = '  +6.46%  '
$r = $ws.Range('D14')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '14.19'
$r.Style = $s
$ws.Range('E14').Value = '  -2.65%  '
$r = $ws.Range('D15')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.260.58'
$r.Style = $s
$ws.Range('E15').Value = '  -4.04%  '
$r = $ws.Range('D16')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '5.25'
$r.Style = $s
$ws.Range('E16').Value = '  -1.53%  '
$r = $ws.Range('D17')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '1.973.62'
$r.Style = $s
$ws.Range('E17').Value = '  -4.00%  '
$ws.Range('E18').Value = '  +5.25%  '
$r = $ws.Range('D19')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '35.425.19'
$r.Style = $s
$ws.Range('E19').Value = '  -2.63%  '
$r = $ws.Range('D20')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '71.25'
$r.Style = $s
$ws.Range('E20').Value = '  -0.86%  '
$r = $ws.Range('D21')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.0₃0839'
$r.Style = $s
$ws.Range('E21').Value = '  -1.64%  '
$r = $ws.Range('D22')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '232.44'
$r.Style = $s
$r = $ws.Range('D23')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '5.13'
$r.Style = $s
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('E24').Value = '  +0.00%  '
$r = $ws.Range('D25')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.56'
$r.Style = $s
$ws.Range('E25').Value = '  +21.14%  '
$ws.Range('E26').Value = '  -1.51%  '
$r = $ws.Range('D27')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '163.31'
$r.Style = $s
$ws.Range('E27').Value = '  +0.36%  '
$r = $ws.Range('D28')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '9.06'
$r.Style = $s
$ws.Range('E28').Value = '  -2.90%  '
$r = $ws.Range('D29')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '19.20'
$r.Style = $s
$ws.Range('E29').Value = '  -4.32%  '
$ws.Range('E30').Value = '  -2.41%  '
$r = $ws.Range('D31')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '4.85'
$r.Style = $s
$ws.Range('E31').Value = '  -3.41%  '
$r = $ws.Range('D32')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '1.13'
$r.Style = $s
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('E33').Value = '  -0.20%  '
$r = $ws.Range('D34')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.0913'
$r.Style = $s
$ws.Range('E34').Value = '  +10.36%  '
$r = $ws.Range('D35')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '4.27'
$r.Style = $s
$ws.Range('E35').Value = '  -4.55%  '
$r = $ws.Range('D36')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.35'
$r.Style = $s
$ws.Range('E36').Value = '  +8.22%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -5.25%  '
$ws.Range('E39').Value = '  +5.31%  '
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('E41').Value = '  +1.19%  '
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('E43').Value = '  -2.19%  '
$r = $ws.Range('D44')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '91.13'
$r.Style = $s
$ws.Range('E44').Value = '  -2.27%  '
$r = $ws.Range('D45')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '1.379.66'
$r.Style = $s
$ws.Range('E45').Value = '  -0.01%  '
$r = $ws.Range('B46')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'InjectiveProtocol'
$r.Style = $s
$r = $ws.Range('C46')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$r.Style = $s
$r = $ws.Range('D46')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '15.89'
$r.Style = $s
$ws.Range('E46').Value = '  +1.86%  '
$r = $ws.Range('B47')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'Cronos'
$r.Style = $s
$r = $ws.Range('C47')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r.Style = $s
$r = $ws.Range('D47')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.0880'
$r.Style = $s
$ws.Range('E47').Value = '  -1.66%  '
$r = $ws.Range('D48')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '7.48'
$r.Style = $s
$ws.Range('E48').Value = '  +3.07%  '
$ws.Range('E49').Value = '  +1.16%  '
$r = $ws.Range('B50')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'FTXToken'
$r.Style = $s
$r = $ws.Range('C50')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$r.Style = $s
$r = $ws.Range('D50')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '3.65'
$r.Style = $s
$ws.Range('E50').Value = '  +9.45%  '
$r = $ws.Range('B51')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'RenderToken'
$r.Style = $s
$r = $ws.Range('C51')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r.Style = $s
$r = $ws.Range('D51')
$s = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.25'
$r.Style = $s
$ws.Range('E51').Value = '  +0.00%  '
